$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.985.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.518.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.519.06"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.422"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.115.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.502.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.986.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  +8.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "437.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.609"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.26%  "
$ws.Range("E24").Value = "  +2.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.642.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.06%  "
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("E31").Value = "  -3.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.168"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.513.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("E36").Value = "  -3.27%  "
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0892"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "171.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -10.55%  "
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.11%  "
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("E50").Value = "  -3.74%  "
$ws.Range("E51").Value = "  -0.04%  "
